$d = $word.ActiveDocument

# --- 1) "Condicions" -> "Condicions de compra" in the equipment-conditions
#        header cell (Table 2, row 1, col 2). Scope the Find/Replace to
#        that single table cell so the other "Condicions ambientals"
#        header elsewhere in the document is left untouched.
$condicionsCell = $d.Tables.Item(2).Cell(1, 2)
$condicionsCell.Range.Find.Execute("Condicions", $true, $false, $false, $false, $false, $true, 0, $false, "Condicions de compra", 1)

# --- 2) Merge "{{body.preu}} " + "€" back into a single run.
$preuCell = $d.Tables.Item(2).Cell(2, 5)
$preuCell.Range.Find.Execute("{{body.preu}} €", $true, $false, $false, $false, $false, $true, 0, $false, "{{body.preu}} €", 1)

# --- 3) Merge "{{body.amplada}} " + "cm" back into a single run.
$ampladaCell = $d.Tables.Item(3).Cell(2, 1)
$ampladaCell.Range.Find.Execute("{{body.amplada}} cm", $true, $false, $false, $false, $false, $true, 0, $false, "{{body.amplada}} cm", 1)

# --- 4) Merge "{{body.alçada}} " + "cm" back into a single run.
$alcadaCell = $d.Tables.Item(3).Cell(2, 2)
$alcadaCell.Range.Find.Execute("{{body.alçada}} cm", $true, $false, $false, $false, $false, $true, 0, $false, "{{body.alçada}} cm", 1)

# --- 5) Merge "{{body.profunditat}} " + "cm" back into a single run.
$profunditatCell = $d.Tables.Item(3).Cell(2, 3)
$profunditatCell.Range.Find.Execute("{{body.profunditat}} cm", $true, $false, $false, $false, $false, $true, 0, $false, "{{body.profunditat}} cm", 1)

# --- 6) Merge "{{body.pes}} " + "Kg" back into a single run.
$pesCell = $d.Tables.Item(3).Cell(2, 4)
$pesCell.Range.Find.Execute("{{body.pes}} Kg", $true, $false, $false, $false, $false, $true, 0, $false, "{{body.pes}} Kg", 1)

# --- 7) Merge "{{body.volum}} " + "L" back into a single run.
$volumCell = $d.Tables.Item(3).Cell(2, 5)
$volumCell.Range.Find.Execute("{{body.volum}} L", $true, $false, $false, $false, $false, $true, 0, $false, "{{body.volum}} L", 1)
